# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> New Value mapping (old values are overwritten regardless of current content)
$updates = @{
    2  = 8020
    3  = 7635
    4  = 111
    9  = 107
    10 = 152
    12 = 690
    13 = 112
    14 = 1215
    15 = 58
    19 = 102
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
